$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "running capri leggings"
$ws.Range("A2").Value = "running capri leggings for women"
$ws.Range("A3").Value = "running capri pants"
$ws.Range("A4").Value = "running capri pants women"
$ws.Range("A5").Value = "running capris"
$ws.Range("A6").Value = "running capris for women"
$ws.Range("A7").Value = "running capris for women with pockets"
$ws.Range("A8").Value = "running capris with side pockets"
$ws.Range("A9").Value = "running capris women"
$ws.Range("A10").Value = "running capris women pocket"
$ws.Range("A11").Value = "running chaffing stick"
$ws.Range("A12").Value = "running chafing"
$ws.Range("A13").Value = "running cloth"
$ws.Range("A14").Value = "running clothes"
$ws.Range("A15").Value = "running clothes for women"
$ws.Range("A16").Value = "running clothes for women cold weather"
$ws.Range("A17").Value = "running clothes for women winter"
$ws.Range("A18").Value = "running clothes reflective"
$ws.Range("A19").Value = "running clothes women"
$ws.Range("A20").Value = "running clothes women winter"
$ws.Range("A21").Value = "running clothing"
$ws.Range("A22").Value = "running clothing women"
$ws.Range("A23").Value = "running cold"
$ws.Range("A24").Value = "running cold gear men"
$ws.Range("A25").Value = "running cold weather gear women"
$ws.Range("A26").Value = "running cold weather pants"
$ws.Range("A27").Value = "running cold weather pants women"
$ws.Range("A28").Value = "running compresion"
$ws.Range("A29").Value = "running compression"
$ws.Range("A30").Value = "running compression calf women"
$ws.Range("A31").Value = "running compression capris women"
$ws.Range("A32").Value = "running compression gear"
$ws.Range("A33").Value = "running compression leg"
$ws.Range("A34").Value = "running compression leggings"
$ws.Range("A35").Value = "running compression leggings women"
$ws.Range("A36").Value = "running compression pants"
$ws.Range("A37").Value = "running compression pants women"
$ws.Range("A38").Value = "running compression shorts"
$ws.Range("A39").Value = "running compression shorts for women"
$ws.Range("A40").Value = "running compression shorts women"
$ws.Range("A41").Value = "running compression tights"
$ws.Range("A42").Value = "running compression tights men"
$ws.Range("A43").Value = "running compression tights women"
$ws.Range("A44").Value = "running compression women"
$ws.Range("A45").Value = "running crop pants"
$ws.Range("A46").Value = "running doesn't suck"
$ws.Range("A47").Value = "running endurance supplements"
$ws.Range("A48").Value = "running equipment for women"
$ws.Range("A49").Value = "running equipment women"
$ws.Range("A50").Value = "running events"
$ws.Range("A51").Value = "running exercise pants women"
$ws.Range("A52").Value = "running faster"
$ws.Range("A53").Value = "running faster training"
$ws.Range("A54").Value = "running gear clothes"
$ws.Range("A55").Value = "running gear cold weather"
$ws.Range("A56").Value = "running gear cold weather women"
$ws.Range("A57").Value = "running gear cold womens"
$ws.Range("A58").Value = "running gear for cold weather"
$ws.Range("A59").Value = "running gear for men cold weather"
$ws.Range("A60").Value = "running gear for winter women"
$ws.Range("A61").Value = "running gear for woman"
$ws.Range("A62").Value = "running gear for women cold"
$ws.Range("A63").Value = "running gear for women cold weather"
$ws.Range("A64").Value = "running gear for women summer"
$ws.Range("A65").Value = "running gear for women winter"
$ws.Range("A66").Value = "running gear gifts for women under 100 dollars"
$ws.Range("A67").Value = "running gear in winter"
$ws.Range("A68").Value = "running gear ladies"
$ws.Range("A69").Value = "running gear marathon"
$ws.Range("A70").Value = "running gear pants"
$ws.Range("A71").Value = "running gear teens"
$ws.Range("A72").Value = "running gear winter"
$ws.Range("A73").Value = "running gear winter women"
$ws.Range("A74").Value = "running gear women"
$ws.Range("A75").Value = "running gears for men"
$ws.Range("A76").Value = "running half tight"
$ws.Range("A77").Value = "running half tights"
$ws.Range("A78").Value = "running half tights men"
$ws.Range("A79").Value = "running half tights women"
$ws.Range("A80").Value = "running hip light"
$ws.Range("A81").Value = "running hip pack"
$ws.Range("A82").Value = "running in cold weather gear"
$ws.Range("A83").Value = "running jacket women reflective"
$ws.Range("A84").Value = "running knee"
$ws.Range("A85").Value = "running knee band"
$ws.Range("A86").Value = "running knee brace for women"
$ws.Range("A87").Value = "running knee brace runners knee"
$ws.Range("A88").Value = "running knee compression"
$ws.Range("A89").Value = "running knee support"
$ws.Range("A90").Value = "running knee support pair"
$ws.Range("A91").Value = "running knee support women"
$ws.Range("A92").Value = "running knee supports"
$ws.Range("A93").Value = "running knees"
$ws.Range("A94").Value = "running legging"
$ws.Range("A95").Value = "running legging women"
$ws.Range("A96").Value = "running leggings"
$ws.Range("A97").Value = "running leggings compression"
$ws.Range("A98").Value = "running leggings compression women"
$ws.Range("A99").Value = "running leggings for women"
$ws.Range("A100").Value = "running leggings for women capri"
